$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("is_active") currently stores =TRUE() formulas that evaluate to
# the numeric boolean 1. The fix replaces them with the literal text value
# "TRUE" while keeping the column's existing "Text" cell formatting.
#
# Typing the bare word TRUE would normally be auto-recognized as the
# boolean literal TRUE again, so a formula that produces the text "TRUE"
# is used first and then converted in place to a plain value. This avoids
# any boolean re-interpretation and keeps the original cell formatting
# untouched (no quote-prefix / number-format side effects).
$range = $ws.Range("E2:E11")
$range.Formula = "=TRIM(""TRUE "")"
$range.Copy()
$range.PasteSpecial(-4163) # xlPasteValues

# Reflect the selection left behind after editing E2:E11.
[void]$range.Select()
